$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean)
$ws.Range("B3").Value = 8497.212576938084
$ws.Range("D3").Value = 576.0166857502905

# Row 4 (std)
$ws.Range("B4").Value = 3628.833582850662
$ws.Range("D4").Value = 545.8823715427411

# Row 5 (min)
$ws.Range("B5").Value = 888.021

# Row 6 (25%)
$ws.Range("B6").Value = 5675.002

# Row 7 (50%)
$ws.Range("B7").Value = 7778.012000000001
$ws.Range("D7").Value = 720

# Row 8 (75%)
$ws.Range("B8").Value = 11531.60025000007
$ws.Range("D8").Value = 1040

# Row 9 (max)
$ws.Range("B9").Value = 19255.31000000004
$ws.Range("D9").Value = 5100

# Row 10 (Total)
$ws.Range("F10").Value = 4466134930.441011

# Row 11 (Residential)
$ws.Range("G11").Value = 0.7143008952600087

# Row 12 (Community)
$ws.Range("F12").Value = 302754370.0310001
$ws.Range("G12").Value = 0.06778889906962673

# Row 13 (IGA)
$ws.Range("G13").Value = 0.2179102056703646
